$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 2 4 '66.287.81'  # D2: '66.323.73' -> '66.287.81'
Set-TextValue 2 5 '  -1.35%  '  # E2: '  -1.42%  ' -> '  -1.35%  '
Set-TextValue 3 4 '3.560.13'  # D3: '3.563.00' -> '3.560.13'
Set-TextValue 4 5 '  -0.05%  '  # E4: '  -0.03%  ' -> '  -0.05%  '
Set-TextValue 5 4 '608.40'  # D5: '609.16' -> '608.40'
Set-TextValue 5 5 '  -0.31%  '  # E5: '  -0.34%  ' -> '  -0.31%  '
Set-TextValue 6 4 '144.50'  # D6: '144.56' -> '144.50'
Set-TextValue 6 5 '  -2.43%  '  # E6: '  -2.54%  ' -> '  -2.43%  '
Set-TextValue 7 4 '3.556.76'  # D7: '3.564.45' -> '3.556.76'
Set-TextValue 7 5 '  +1.29%  '  # E7: '  +1.57%  ' -> '  +1.29%  '
Set-TextValue 8 5 '  -0.01%  '  # E8: '  -0.03%  ' -> '  -0.01%  '
Set-TextValue 9 5 '  +0.38%  '  # E9: '  +0.29%  ' -> '  +0.38%  '
Set-TextValue 11 4 '8.08'  # D11: '8.09' -> '8.08'
Set-TextValue 11 5 '  +0.23%  '  # E11: '  +0.47%  ' -> '  +0.23%  '
Set-TextValue 12 5 '  -2.82%  '  # E12: '  -2.69%  ' -> '  -2.82%  '
Set-TextValue 13 4 '4.160.66'  # D13: '4.166.40' -> '4.160.66'
Set-TextValue 13 5 '  +1.23%  '  # E13: '  +1.30%  ' -> '  +1.23%  '
Set-TextValue 14 5 '  -4.04%  '  # E14: '  -3.99%  ' -> '  -4.04%  '
Set-TextValue 15 4 '30.16'  # D15: '30.19' -> '30.16'
Set-TextValue 15 5 '  -4.46%  '  # E15: '  -4.34%  ' -> '  -4.46%  '
Set-TextValue 16 4 '3.561.90'  # D16: '3.558.12' -> '3.561.90'
Set-TextValue 16 5 '  +1.30%  '  # E16: '  +1.12%  ' -> '  +1.30%  '
Set-TextValue 17 4 '66.364.84'  # D17: '66.397.62' -> '66.364.84'
Set-TextValue 17 5 '  -1.37%  '  # E17: '  -1.40%  ' -> '  -1.37%  '
Set-TextValue 18 5 '  -1.03%  '  # E18: '  -1.04%  ' -> '  -1.03%  '
Set-TextValue 19 4 '11.32'  # D19: '11.36' -> '11.32'
Set-TextValue 19 5 '  +4.04%  '  # E19: '  +4.14%  ' -> '  +4.04%  '
Set-TextValue 20 4 '6.21'  # D20: '6.22' -> '6.21'
Set-TextValue 20 5 '  -2.40%  '  # E20: '  -2.24%  ' -> '  -2.40%  '
Set-TextValue 21 4 '14.95'  # D21: '14.98' -> '14.95'
Set-TextValue 21 5 '  -3.13%  '  # E21: '  -2.92%  ' -> '  -3.13%  '
Set-TextValue 22 4 '428.97'  # D22: '429.41' -> '428.97'
Set-TextValue 22 5 '  -1.75%  '  # E22: '  -1.65%  ' -> '  -1.75%  '
Set-TextValue 23 5 '  -0.80%  '  # E23: '  -0.71%  ' -> '  -0.80%  '
Set-TextValue 24 4 '78.92'  # D24: '78.96' -> '78.92'
Set-TextValue 25 4 '3.697.36'  # D25: '3.700.25' -> '3.697.36'
Set-TextValue 25 5 '  +1.14%  '  # E25: '  +1.16%  ' -> '  +1.14%  '
Set-TextValue 26 5 '  -0.04%  '  # E26: '  -0.03%  ' -> '  -0.04%  '
Set-TextValue 27 5 '  +2.42%  '  # E27: '  +2.73%  ' -> '  +2.42%  '
Set-TextValue 28 4 '8.10'  # D28: '8.12' -> '8.10'
Set-TextValue 28 5 '  -2.18%  '  # E28: '  -1.93%  ' -> '  -2.18%  '
Set-TextValue 29 5 '  -6.52%  '  # E29: '  -6.45%  ' -> '  -6.52%  '
Set-TextValue 30 5 '  -1.27%  '  # E30: '  -1.33%  ' -> '  -1.27%  '
Set-TextValue 31 4 '0.999'  # D31: '1.00' -> '0.999'
Set-TextValue 31 5 '  -0.09%  '  # E31: '  +0.08%  ' -> '  -0.09%  '
Set-TextValue 32 5 '  -5.72%  '  # E32: '  -5.52%  ' -> '  -5.72%  '
Set-TextValue 33 5 '  -4.15%  '  # E33: '  -3.91%  ' -> '  -4.15%  '
Set-TextValue 34 4 '25.48'  # D34: '25.45' -> '25.48'
Set-TextValue 34 5 '  -0.42%  '  # E34: '  -0.51%  ' -> '  -0.42%  '
Set-TextValue 35 4 '3.548.83'  # D35: '3.552.33' -> '3.548.83'
Set-TextValue 35 5 '  +1.18%  '  # E35: '  +1.22%  ' -> '  +1.18%  '
Set-TextValue 36 5 '  -0.04%  '  # E36: '  -0.03%  ' -> '  -0.04%  '
Set-TextValue 37 5 '  -3.29%  '  # E37: '  -3.03%  ' -> '  -3.29%  '
Set-TextValue 38 2 'NEARProtocol'  # B38: 'Aptos' -> 'NEARProtocol'
Set-TextValue 38 3 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'  # C38: 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt' -> 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 38 4 '5.64'  # D38: '7.84' -> '5.64'
Set-TextValue 38 5 '  -5.49%  '  # E38: '  -2.38%  ' -> '  -5.49%  '
Set-TextValue 39 2 'Aptos'  # B39: 'NEARProtocol' -> 'Aptos'
Set-TextValue 39 3 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'  # C39: 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near' -> 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 39 4 '7.83'  # D39: '5.65' -> '7.83'
Set-TextValue 39 5 '  -2.67%  '  # E39: '  -5.39%  ' -> '  -2.67%  '
Set-TextValue 40 4 '1.00'  # D40: '0.999' -> '1.00'
Set-TextValue 40 5 '  -0.05%  '  # E40: '  +0.02%  ' -> '  -0.05%  '
Set-TextValue 41 4 '174.88'  # D41: '174.74' -> '174.88'
Set-TextValue 41 5 '  -0.68%  '  # E41: '  -0.91%  ' -> '  -0.68%  '
Set-TextValue 42 4 '0.0859'  # D42: '0.0861' -> '0.0859'
Set-TextValue 42 5 '  -4.72%  '  # E42: '  -4.48%  ' -> '  -4.72%  '
Set-TextValue 43 4 '5.27'  # D43: '5.28' -> '5.27'
Set-TextValue 43 5 '  -2.58%  '  # E43: '  -2.34%  ' -> '  -2.58%  '
Set-TextValue 44 4 '0.895'  # D44: '0.896' -> '0.895'
Set-TextValue 44 5 '  -0.12%  '  # E44: '  +0.03%  ' -> '  -0.12%  '
Set-TextValue 45 5 '  -7.22%  '  # E45: '  -6.97%  ' -> '  -7.22%  '
Set-TextValue 46 4 '45.67'  # D46: '45.70' -> '45.67'
Set-TextValue 46 5 '  -1.54%  '  # E46: '  -1.44%  ' -> '  -1.54%  '
Set-TextValue 47 5 '  -1.64%  '  # E47: '  -1.41%  ' -> '  -1.64%  '
Set-TextValue 48 4 '26.03'  # D48: '26.07' -> '26.03'
Set-TextValue 48 5 '  -8.95%  '  # E48: '  -7.60%  ' -> '  -8.95%  '
Set-TextValue 49 4 '2.39'  # D49: '2.40' -> '2.39'
Set-TextValue 49 5 '  -2.48%  '  # E49: '  -2.30%  ' -> '  -2.48%  '
Set-TextValue 50 4 '7.14'  # D50: '7.15' -> '7.14'
Set-TextValue 50 5 '  -4.49%  '  # E50: '  -4.37%  ' -> '  -4.49%  '
Set-TextValue 51 4 '23.01'  # D51: '23.05' -> '23.01'
Set-TextValue 51 5 '  +6.33%  '  # E51: '  +6.52%  ' -> '  +6.33%  '
